# Daily attendance processing - 2025-10-24 14:48:21
# Normalize the "Recorded By" (column G) comma-separated author lists for the
# records still holding these two specific legacy orderings by rotating the
# list right by one (moving the last entry to the front).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 160
$changed = 0

# Exact legacy values that need to be re-ordered, mapped to their new value.
$map = @{
    'System, dnasr281@gmail.com'          = 'dnasr281@gmail.com, System'
    'System, backup@backdoor.com, system' = 'system, System, backup@backdoor.com'
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Value2

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if ($map.ContainsKey($text)) {
        $cell.Value2 = $map[$text]
        $changed++
    }
}

Write-Host "Rotated Recorded By values in $changed cells"
